# ---------------------------------------------------------------------
# Weekly NYPD CompStat source update (83rd Precinct): new week of crime
# data collected (Volume 30, Number 15 -- week covering 4/10/2023 through
# 4/16/2023). Updates the report header and every Week-to-Date / 28-Day /
# Year-to-Date / 2-Year count + percent-change figure in the crime table.
# ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Report header text: bump the issue number and the week-covering
#    date range. Characters() edits only the trailing "14" -> "15" in
#    place; the week-range cell is replaced outright with the new text.
# -----------------------------------------------------------------
$issueCell = $ws.Range("A8")
$issueText = $issueCell.Text
$issueCell.Characters($issueText.Length - 1, 2).Text = "15"

$ws.Range("C9").Value = "Report Covering the Week  4/10/2023  Through  4/16/2023"

# -----------------------------------------------------------------
# 2) Helper functions to change a cell's stored type (number <-> text)
#    while keeping the destination cell's number format in line with
#    its new content. Excel stores count columns either as a number
#    (style carrying numFmt "#,##0") or, when the count is zero and
#    the week-over-week change is undefined, as the literal text
#    placeholders "0" / "***.*" (General-format style). Setting
#    .Value first fixes the stored type, then copying +
#    PasteSpecial(xlPasteFormats) from a same-shaped template cell
#    re-applies the correct number format/style without touching the
#    value just written.
# -----------------------------------------------------------------
$xlPasteFormats = -4122

function Set-TextPlaceholder($ref, $text) {
    $dst = $ws.Range($ref)
    $dst.Value = "'" + $text
    $ws.Range("C14").Copy()
    $dst.PasteSpecial($xlPasteFormats)
}

function Set-NumberValue($ref, $value, $templateRef) {
    $dst = $ws.Range($ref)
    $dst.Value = $value
    $ws.Range($templateRef).Copy()
    $dst.PasteSpecial($xlPasteFormats)
}

# Cells that flip from a count (number) to the "no prior-year activity"
# text placeholders ("0" / "***.*"), matching the style (14) + shared
# string used elsewhere in the sheet for the same placeholders.
Set-TextPlaceholder "D15" "0"
Set-TextPlaceholder "E15" "***.*"
Set-TextPlaceholder "C22" "0"
Set-TextPlaceholder "D22" "0"
Set-TextPlaceholder "E22" "***.*"
Set-TextPlaceholder "G23" "0"
Set-TextPlaceholder "H23" "***.*"
Set-TextPlaceholder "D26" "0"
Set-TextPlaceholder "E26" "***.*"
Set-TextPlaceholder "D30" "0"
Set-TextPlaceholder "E30" "***.*"

# Cells that flip the other way: from the text placeholders back to an
# actual numeric count / percent-change value.
Set-NumberValue "C18" 3 "D16"
Set-NumberValue "D27" 3 "D16"
Set-NumberValue "E27" 0 "H16"

# -----------------------------------------------------------------
# 3) Bulk numeric updates: cells whose style (number format) stays the
#    same -- only the stored count / percent-change value moved this
#    week.
# -----------------------------------------------------------------
$numericUpdates = @(
    @("N14", -92.857142857142),
    @("G15", 3),
    @("L15", -28.571428571428),
    @("N15", -78.260869565217),
    @("C16", 9),
    @("D16", 6),
    @("E16", 50),
    @("F16", 29),
    @("G16", 14),
    @("H16", 107.142857142857),
    @("I16", 86),
    @("J16", 76),
    @("K16", 13.157894736842),
    @("L16", 21.126760563380),
    @("M16", -26.495726495726),
    @("N16", -79.126213592233),
    @("C17", 4),
    @("D17", 6),
    @("E17", -33.333333333333),
    @("F17", 21),
    @("G17", 29),
    @("H17", -27.586206896551),
    @("I17", 110),
    @("J17", 97),
    @("K17", 13.402061855670),
    @("L17", 30.952380952381),
    @("M17", 41.025641025641),
    @("N17", -54.918032786885),
    @("D18", 3),
    @("E18", 0),
    @("F18", 9),
    @("H18", -57.142857142857),
    @("I18", 64),
    @("J18", 83),
    @("K18", -22.891566265060),
    @("L18", -20.987654320987),
    @("M18", -26.436781609195),
    @("N18", -83.715012722646),
    @("D19", 11),
    @("E19", -9.090909090909),
    @("F19", 41),
    @("G19", 53),
    @("H19", -22.641509433962),
    @("I19", 199),
    @("J19", 162),
    @("K19", 22.839506172839),
    @("L19", 70.085470085470),
    @("M19", 136.904761904762),
    @("N19", 34.459459459459),
    @("C20", 5),
    @("D20", 5),
    @("F20", 14),
    @("H20", -12.5),
    @("I20", 49),
    @("J20", 60),
    @("K20", -18.333333333333),
    @("L20", 36.111111111111),
    @("M20", 13.953488372093),
    @("N20", -83.774834437086),
    @("C21", 31),
    @("D21", 31),
    @("E21", 0),
    @("F21", 114),
    @("G21", 136),
    @("H21", -16.176470588235),
    @("I21", 514),
    @("J21", 491),
    @("K21", 4.684317718940),
    @("L21", 29.471032745591),
    @("M21", 22.966507177033),
    @("N21", -66.536458333333),
    @("F22", 2),
    @("H22", 100),
    @("M22", 11.111111111111),
    @("L23", -27.272727272727),
    @("M23", 60),
    @("C24", 21),
    @("D24", 19),
    @("E24", 10.526315789473),
    @("F24", 78),
    @("G24", 66),
    @("H24", 18.181818181818),
    @("I24", 263),
    @("J24", 261),
    @("K24", 0.766283524904),
    @("L24", 5.622489959839),
    @("M24", 35.567010309278),
    @("C25", 14),
    @("D25", 12),
    @("E25", 16.666666666666),
    @("F25", 46),
    @("G25", 51),
    @("H25", -9.803921568627),
    @("I25", 177),
    @("J25", 184),
    @("K25", -3.804347826086),
    @("L25", 43.902439024390),
    @("M25", -13.235294117647),
    @("G26", 4),
    @("H26", 0),
    @("L26", -23.076923076923),
    @("C27", 3),
    @("F27", 5),
    @("G27", 8),
    @("H27", -37.5),
    @("I27", 15),
    @("J27", 17),
    @("K27", -11.764705882352),
    @("L27", -6.25),
    @("M28", -77.777777777777),
    @("N28", -96.774193548387),
    @("M29", -75),
    @("N29", -96.721311475409)
)

foreach ($pair in $numericUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

Write-Host "Weekly crime data refreshed: $($numericUpdates.Count) figures updated, plus 14 placeholder/number style swaps and the report header."
